$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Coin (B), Link (C), Price (D), Volume(1h) (E)
$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '63.185.68', '  -5.98%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '3.467.78', '  -2.90%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  -0.20%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '387.96', '  -6.50%  '),
    @(6, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '122.42', '  -4.74%  '),
    @(7, 'LidoStakedEther', 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth', '3.460.71', '  -2.90%  '),
    @(8, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.583', '  -9.89%  '),
    @(9, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.999', '  +0.06%  '),
    @(10, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.667', '  -12.77%  '),
    @(11, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.143', '  -17.44%  '),
    @(12, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000317', '  -2.98%  '),
    @(13, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '38.50', '  -8.30%  '),
    @(14, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '3.992.11', '  -3.73%  '),
    @(15, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '9.13', '  -6.79%  '),
    @(16, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.136', '  -3.20%  '),
    @(17, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '3.455.33', '  -3.08%  '),
    @(18, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '12.49', '  +2.10%  '),
    @(19, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '18.46', '  -8.95%  '),
    @(20, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '63.171.65', '  -5.91%  '),
    @(21, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '1.01', '  -10.51%  '),
    @(22, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '390.32', '  -13.71%  '),
    @(23, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '13.72', '  +3.47%  '),
    @(24, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '79.84', '  -9.44%  '),
    @(25, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '2.83', '  -9.58%  '),
    @(26, 'LEO', 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo', '5.21', '  +6.76%  '),
    @(27, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '32.79', '  -5.20%  '),
    @(28, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '2.94', '  -12.23%  '),
    @(29, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '8.79', '  -12.81%  '),
    @(30, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '11.76', '  -4.02%  '),
    @(31, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.60', '  -6.17%  '),
    @(32, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.108', '  -6.93%  '),
    @(33, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '6.70', '  -8.43%  '),
    @(34, 'Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.152', '  -5.61%  '),
    @(35, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.00', '  +0.16%  '),
    @(36, 'InjectiveProtocol', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj', '36.45', '  -10.16%  '),
    @(37, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '53.65', '  -5.19%  '),
    @(38, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0434', '  -11.61%  '),
    @(39, 'FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '0.995', '  -0.47%  '),
    @(40, 'PEPE', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', '0.0₃0630', '  -12.69%  '),
    @(41, 'ThetaToken', 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta', '2.64', '  +13.92%  '),
    @(42, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.130', '  -10.67%  '),
    @(43, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '138.39', '  -7.18%  '),
    @(44, 'ApeXProtocol', 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex', '2.99', '  +11.70%  '),
    @(45, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.94', '  -1.34%  '),
    @(46, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '3.05', '  -5.86%  '),
    @(47, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '24.52', '  +14.69%  '),
    @(48, 'WEMIXToken', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '2.45', '  -10.69%  '),
    @(49, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '4.00', '  -6.66%  '),
    @(50, 'Stacks', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', '2.64', '  -13.00%  '),
    @(51, 'TheGraph', 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt', '0.274', '  -11.53%  '),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "'" + $row[3]
    $ws.Cells.Item($r, 5).Value = "'" + $row[4]
}

$wb.Save()